$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 49, shifting existing rows 49.. down to 50..
$ws.Rows("49:49").Insert()

# After the insert, the old row-49 data now lives in row 50. Copy that
# entire row's values/format into the new (currently blank) row 49, then
# change only the Fecha (date) column to the new value.
$src = $ws.Range("A50:T50")
$dst = $ws.Range("A49:T49")
$src.Copy()
$dst.PasteSpecial(-4104)  # xlPasteAll

$ws.Range("D49").Value = 44662
